$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investment_Cost")

# Rename electrolyzer object type labels (reordered naming convention)
$ws.Range("A3").Value = "AEC_Electrolyzer"
$ws.Range("A4").Value = "PEM_Electrolyzer"
$ws.Range("A5").Value = "SOEC_Electrolyzer"

# Update the active selection as in the saved workbook
$ws.Range("A6").Select()
